$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four "Resolving-Mac" sending-cluster rows (old rows 14-17);
# Excel shifts subsequent rows up automatically on each Delete().
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Delete()

# Recomputed TPM-based NATMI metrics for the 12 remaining Il34-Csf1r rows.
# Row 2: ECs -> ECs
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.471704
$ws.Range("H2").Value = 4.415112
$ws.Range("I2").Value = 0.08657490103749592
$ws.Range("J2").Value = 0.0865749010374959
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3667156666666667
$ws.Range("N2").Value = 1.100147
$ws.Range("O2").Value = 0.001893484040582144
$ws.Range("P2").Value = 0.001893484040582144
$ws.Range("Q2").Value = 0.5396969134959999
$ws.Range("R2").Value = 4.857272221463999
$ws.Range("S2").Value = 0.000163928193429477
$ws.Range("T2").Value = 0.000163928193429477

# Row 3: ECs -> FAPs
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.471704
$ws.Range("H3").Value = 4.415112
$ws.Range("I3").Value = 0.08657490103749592
$ws.Range("J3").Value = 0.0865749010374959
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.853217
$ws.Range("N3").Value = 5.559651000000001
$ws.Range("O3").Value = 0.009568821657202681
$ws.Range("P3").Value = 0.00956882165720268
$ws.Range("Q3").Value = 2.727386871768
$ws.Range("R3").Value = 24.546481845912
$ws.Range("S3").Value = 0.0008284197880177698
$ws.Range("T3").Value = 0.0008284197880177696

# Row 4: ECs -> MuSCs
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.471704
$ws.Range("H4").Value = 4.415112
$ws.Range("I4").Value = 0.08657490103749592
$ws.Range("J4").Value = 0.0865749010374959
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.599526666666667
$ws.Range("N4").Value = 4.79858
$ws.Range("O4").Value = 0.008258927804608534
$ws.Range("P4").Value = 0.008258927804608534
$ws.Range("Q4").Value = 2.35402979344
$ws.Range("R4").Value = 21.18626814096
$ws.Range("S4").Value = 0.0007150158573598072
$ws.Range("T4").Value = 0.0007150158573598071

# Row 5: ECs -> Resolving-Mac
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.471704
$ws.Range("H5").Value = 4.415112
$ws.Range("I5").Value = 0.08657490103749592
$ws.Range("J5").Value = 0.0865749010374959
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 189.8529766666667
$ws.Range("N5").Value = 569.55893
$ws.Range("O5").Value = 0.9802787664976067
$ws.Range("P5").Value = 0.9802787664976066
$ws.Range("Q5").Value = 279.40738517224
$ws.Range("R5").Value = 2514.66646655016
$ws.Range("S5").Value = 0.08486753719868886
$ws.Range("T5").Value = 0.08486753719868885

# Row 6: FAPs -> ECs
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.873786
$ws.Range("H6").Value = 14.621358
$ws.Range("I6").Value = 0.2867067974456365
$ws.Range("J6").Value = 0.2867067974456365
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3667156666666667
$ws.Range("N6").Value = 1.100147
$ws.Range("O6").Value = 0.001893484040582144
$ws.Range("P6").Value = 0.001893484040582144
$ws.Range("Q6").Value = 1.787293682180667
$ws.Range("R6").Value = 16.085643139626
$ws.Range("S6").Value = 0.0005428747452897302
$ws.Range("T6").Value = 0.0005428747452897302

# Row 7: FAPs -> FAPs
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.873786
$ws.Range("H7").Value = 14.621358
$ws.Range("I7").Value = 0.2867067974456365
$ws.Range("J7").Value = 0.2867067974456365
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.853217
$ws.Range("N7").Value = 5.559651000000001
$ws.Range("O7").Value = 0.009568821657202681
$ws.Range("P7").Value = 0.00956882165720268
$ws.Range("Q7").Value = 9.032183069562
$ws.Range("R7").Value = 81.28964762605801
$ws.Range("S7").Value = 0.002743446212665029
$ws.Range("T7").Value = 0.002743446212665029

# Row 8: FAPs -> MuSCs
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.873786
$ws.Range("H8").Value = 14.621358
$ws.Range("I8").Value = 0.2867067974456365
$ws.Range("J8").Value = 0.2867067974456365
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.599526666666667
$ws.Range("N8").Value = 4.79858
$ws.Range("O8").Value = 0.008258927804608534
$ws.Range("P8").Value = 0.008258927804608534
$ws.Range("Q8").Value = 7.795750674626667
$ws.Range("R8").Value = 70.16175607164001
$ws.Range("S8").Value = 0.002367890741194035
$ws.Range("T8").Value = 0.002367890741194035

# Row 9: FAPs -> Resolving-Mac
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.873786
$ws.Range("H9").Value = 14.621358
$ws.Range("I9").Value = 0.2867067974456365
$ws.Range("J9").Value = 0.2867067974456365
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 189.8529766666667
$ws.Range("N9").Value = 569.55893
$ws.Range("O9").Value = 0.9802787664976067
$ws.Range("P9").Value = 0.9802787664976066
$ws.Range("Q9").Value = 925.3027797363267
$ws.Range("R9").Value = 8327.72501762694
$ws.Range("S9").Value = 0.2810525857464878
$ws.Range("T9").Value = 0.2810525857464877

# Row 10: MuSCs -> ECs
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 10.65370933333333
$ws.Range("H10").Value = 31.961128
$ws.Range("I10").Value = 0.6267183015168676
$ws.Range("J10").Value = 0.6267183015168675
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3667156666666667
$ws.Range("N10").Value = 1.100147
$ws.Range("O10").Value = 0.001893484040582144
$ws.Range("P10").Value = 0.001893484040582144
$ws.Range("Q10").Value = 3.906882120646222
$ws.Range("R10").Value = 35.161939085816
$ws.Range("S10").Value = 0.001186681101862937
$ws.Range("T10").Value = 0.001186681101862936

# Row 11: MuSCs -> FAPs
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 10.65370933333333
$ws.Range("H11").Value = 31.961128
$ws.Range("I11").Value = 0.6267183015168676
$ws.Range("J11").Value = 0.6267183015168675
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.853217
$ws.Range("N11").Value = 5.559651000000001
$ws.Range("O11").Value = 0.009568821657202681
$ws.Range("P11").Value = 0.00956882165720268
$ws.Range("Q11").Value = 19.743635249592
$ws.Range("R11").Value = 177.692717246328
$ws.Range("S11").Value = 0.005996955656519883
$ws.Range("T11").Value = 0.00599695565651988

# Row 12: MuSCs -> MuSCs
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 10.65370933333333
$ws.Range("H12").Value = 31.961128
$ws.Range("I12").Value = 0.6267183015168676
$ws.Range("J12").Value = 0.6267183015168675
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.599526666666667
$ws.Range("N12").Value = 4.79858
$ws.Range("O12").Value = 0.008258927804608534
$ws.Range("P12").Value = 0.008258927804608534
$ws.Range("Q12").Value = 17.04089217758222
$ws.Range("R12").Value = 153.36802959824
$ws.Range("S12").Value = 0.005176021206054692
$ws.Range("T12").Value = 0.005176021206054691

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 10.65370933333333
$ws.Range("H13").Value = 31.961128
$ws.Range("I13").Value = 0.6267183015168676
$ws.Range("J13").Value = 0.6267183015168675
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 189.8529766666667
$ws.Range("N13").Value = 569.55893
$ws.Range("O13").Value = 0.9802787664976067
$ws.Range("P13").Value = 0.9802787664976066
$ws.Range("Q13").Value = 2022.638429474782
$ws.Range("R13").Value = 18203.74586527304
$ws.Range("S13").Value = 0.6143586435524301
$ws.Range("T13").Value = 0.61435864355243

